$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04215534119371403
$ws.Range("D2").Value = 0.1361288253571666
$ws.Range("G2").Value = 0.1217136106832186
$ws.Range("H2").Value = 0.9740000000000001
